# This script applies the roster/odds refresh captured in the commit:
#   "Atualizado por script em 05-11-2023 14:45"
# It touches only the betting-odds sheet that already exists in the workbook
# (Sheet1). The source scraper re-ran and, for a handful of already-scraped
# fixtures, rows got written back in a different order (rows within the same
# matchday group got shuffled) and two brand-new fixtures were appended at
# the end (rows 117-118), extending the used range from A1:V116 to A1:V118.
#
# Strategy: for every affected existing row we overwrite columns F:V (home
# team .. match url) with their new values; columns A:E (index/country/
# tournament/season/date) are untouched. For the two new rows we first copy
# the last existing row's formatting down (so the Indice column keeps its
# bordered/bold style and the date column keeps its date number format),
# then overwrite every cell with the real values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").Value = "Puchov"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "Myjava"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2.02
$ws.Range("K3").Value = "28/07/2023 03:42"
$ws.Range("L3").Value = 1.93
$ws.Range("M3").Value = "29/07/2023 16:29"
$ws.Range("N3").Value = 3.36
$ws.Range("O3").Value = "28/07/2023 03:42"
$ws.Range("P3").Value = 3.69
$ws.Range("Q3").Value = "29/07/2023 16:30"
$ws.Range("R3").Value = 3.1
$ws.Range("S3").Value = "28/07/2023 03:42"
$ws.Range("T3").Value = 3.55
$ws.Range("U3").Value = "29/07/2023 16:30"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/slovakia/2-liga/msk-puchov-myjava/2gtPdPJs/"

# Row 4
$ws.Range("F4").Value = "Spisska Nova Ves"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = "Presov"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4.59
$ws.Range("K4").Value = "28/07/2023 03:42"
$ws.Range("L4").Value = 5.61
$ws.Range("M4").Value = "29/07/2023 15:53"
$ws.Range("N4").Value = 3.99
$ws.Range("O4").Value = "28/07/2023 03:42"
$ws.Range("P4").Value = 4.81
$ws.Range("Q4").Value = "29/07/2023 15:53"
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = "28/07/2023 03:42"
$ws.Range("T4").Value = 1.46
$ws.Range("U4").Value = "29/07/2023 15:53"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-presov/tnW4iUs2/"

# Row 5
$ws.Range("F5").Value = "FK Humenne"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "Komarno"
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2.83
$ws.Range("K5").Value = "28/07/2023 03:42"
$ws.Range("L5").Value = 2.97
$ws.Range("M5").Value = "29/07/2023 16:11"
$ws.Range("N5").Value = 3.09
$ws.Range("O5").Value = "28/07/2023 03:42"
$ws.Range("P5").Value = 3.19
$ws.Range("Q5").Value = "29/07/2023 16:11"
$ws.Range("R5").Value = 2.28
$ws.Range("S5").Value = "28/07/2023 03:42"
$ws.Range("T5").Value = 2.37
$ws.Range("U5").Value = "29/07/2023 16:11"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-komarno/dfEiFXZ8/"

# Row 92
$ws.Range("F92").Value = "FK Humenne"
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = "Malzenice"
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = 1.53
$ws.Range("K92").Value = "13/10/2023 02:13"
$ws.Range("L92").Value = 1.39
$ws.Range("M92").Value = "14/10/2023 14:39"
$ws.Range("N92").Value = 3.95
$ws.Range("O92").Value = "13/10/2023 02:13"
$ws.Range("P92").Value = 4.47
$ws.Range("Q92").Value = "14/10/2023 14:39"
$ws.Range("R92").Value = 4.82
$ws.Range("S92").Value = "13/10/2023 02:13"
$ws.Range("T92").Value = 7.94
$ws.Range("U92").Value = "14/10/2023 14:39"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-malzenice/z3np52Ui/"

# Row 93
$ws.Range("F93").Value = "D. Kubin"
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = "Spisska Nova Ves"
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 2
$ws.Range("K93").Value = "13/10/2023 02:13"
$ws.Range("L93").Value = 2.45
$ws.Range("M93").Value = "14/10/2023 14:48"
$ws.Range("N93").Value = 3.36
$ws.Range("O93").Value = "13/10/2023 02:13"
$ws.Range("P93").Value = 3.49
$ws.Range("Q93").Value = "14/10/2023 14:51"
$ws.Range("R93").Value = 3.15
$ws.Range("S93").Value = "13/10/2023 02:13"
$ws.Range("T93").Value = 2.65
$ws.Range("U93").Value = "14/10/2023 14:48"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-spisska-nova-ves/IRyk4Mqb/"

# Row 98
$ws.Range("F98").Value = "Spisska Nova Ves"
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = "FK Humenne"
$ws.Range("I98").Value = 3
$ws.Range("J98").Value = 3.02
$ws.Range("K98").Value = "20/10/2023 01:42"
$ws.Range("L98").Value = 3.43
$ws.Range("M98").Value = "21/10/2023 14:27"
$ws.Range("N98").Value = 3.19
$ws.Range("O98").Value = "20/10/2023 01:42"
$ws.Range("P98").Value = 3.47
$ws.Range("Q98").Value = "21/10/2023 14:27"
$ws.Range("R98").Value = 2.13
$ws.Range("S98").Value = "20/10/2023 01:42"
$ws.Range("T98").Value = 2.04
$ws.Range("U98").Value = "21/10/2023 14:27"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/slovakia/2-liga/spisska-nova-ves-fk-humenne/fTpYgxMj/"

# Row 99
$ws.Range("F99").Value = "D. Kubin"
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = "L. Mikulas"
$ws.Range("I99").Value = 2
$ws.Range("J99").Value = 6.42
$ws.Range("K99").Value = "20/10/2023 01:53"
$ws.Range("L99").Value = 4.12
$ws.Range("M99").Value = "21/10/2023 14:29"
$ws.Range("N99").Value = 4.84
$ws.Range("O99").Value = "20/10/2023 01:53"
$ws.Range("P99").Value = 3.58
$ws.Range("Q99").Value = "21/10/2023 14:29"
$ws.Range("R99").Value = 1.33
$ws.Range("S99").Value = "20/10/2023 01:53"
$ws.Range("T99").Value = 1.82
$ws.Range("U99").Value = "21/10/2023 14:29"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-l-mikulas/tQtUfI6p/"

# Row 100
$ws.Range("F100").Value = "Komarno"
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = "Trebisov"
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1.18
$ws.Range("K100").Value = "20/10/2023 01:42"
$ws.Range("L100").Value = 1.28
$ws.Range("M100").Value = "21/10/2023 14:15"
$ws.Range("N100").Value = 6.03
$ws.Range("O100").Value = "20/10/2023 01:42"
$ws.Range("P100").Value = 5.64
$ws.Range("Q100").Value = "21/10/2023 14:17"
$ws.Range("R100").Value = 9.76
$ws.Range("S100").Value = "20/10/2023 01:42"
$ws.Range("T100").Value = 9.09
$ws.Range("U100").Value = "21/10/2023 14:15"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/slovakia/2-liga/komarno-trebisov/pfeb0vqN/"

# Row 101
$ws.Range("F101").Value = "Malzenice"
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = "Myjava"
$ws.Range("I101").Value = 1
$ws.Range("J101").Value = 2.59
$ws.Range("K101").Value = "20/10/2023 01:42"
$ws.Range("L101").Value = 3.04
$ws.Range("M101").Value = "21/10/2023 14:28"
$ws.Range("N101").Value = 3.23
$ws.Range("O101").Value = "20/10/2023 01:42"
$ws.Range("P101").Value = 3.14
$ws.Range("Q101").Value = "21/10/2023 14:28"
$ws.Range("R101").Value = 2.4
$ws.Range("S101").Value = "20/10/2023 01:42"
$ws.Range("T101").Value = 2.36
$ws.Range("U101").Value = "21/10/2023 14:23"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-myjava/4bSOGaT3/"

# Row 102
$ws.Range("F102").Value = "Presov"
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = "Petrzalka"
$ws.Range("I102").Value = 1
$ws.Range("J102").Value = 1.56
$ws.Range("K102").Value = "20/10/2023 01:42"
$ws.Range("L102").Value = 1.85
$ws.Range("M102").Value = "21/10/2023 14:19"
$ws.Range("N102").Value = 3.84
$ws.Range("O102").Value = "20/10/2023 01:42"
$ws.Range("P102").Value = 3.74
$ws.Range("Q102").Value = "21/10/2023 14:19"
$ws.Range("R102").Value = 4.69
$ws.Range("S102").Value = "20/10/2023 01:42"
$ws.Range("T102").Value = 3.82
$ws.Range("U102").Value = "21/10/2023 14:19"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/slovakia/2-liga/presov-petrzalka/Opf2abbT/"

# Row 104
$ws.Range("F104").Value = "D. Kubin"
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = "Presov"
$ws.Range("I104").Value = 3
$ws.Range("J104").Value = 5.32
$ws.Range("K104").Value = "27/10/2023 02:42"
$ws.Range("L104").Value = 8.949999999999999
$ws.Range("M104").Value = "28/10/2023 14:28"
$ws.Range("N104").Value = 3.95
$ws.Range("O104").Value = "27/10/2023 02:42"
$ws.Range("P104").Value = 5.83
$ws.Range("Q104").Value = "28/10/2023 14:28"
$ws.Range("R104").Value = 1.48
$ws.Range("S104").Value = "27/10/2023 02:42"
$ws.Range("T104").Value = 1.26
$ws.Range("U104").Value = "28/10/2023 14:28"
$ws.Range("V104").Value = "https://www.betexplorer.com/football/slovakia/2-liga/d-kubin-presov/ILYvDHSS/"

# Row 107
$ws.Range("F107").Value = "FK Humenne"
$ws.Range("G107").Value = 4
$ws.Range("H107").Value = "L. Mikulas"
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 2.3
$ws.Range("K107").Value = "27/10/2023 02:42"
$ws.Range("L107").Value = 1.96
$ws.Range("M107").Value = "28/10/2023 14:21"
$ws.Range("N107").Value = 3.18
$ws.Range("O107").Value = "27/10/2023 02:42"
$ws.Range("P107").Value = 3.59
$ws.Range("Q107").Value = "28/10/2023 14:21"
$ws.Range("R107").Value = 2.74
$ws.Range("S107").Value = "27/10/2023 02:42"
$ws.Range("T107").Value = 3.52
$ws.Range("U107").Value = "28/10/2023 14:21"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/slovakia/2-liga/fk-humenne-l-mikulas/x6FB7ejj/"

# Row 108
$ws.Range("F108").Value = "Trebisov"
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = "Zilina B"
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 2.77
$ws.Range("K108").Value = "27/10/2023 02:42"
$ws.Range("L108").Value = 2.31
$ws.Range("M108").Value = "28/10/2023 14:15"
$ws.Range("N108").Value = 3.43
$ws.Range("O108").Value = "27/10/2023 02:42"
$ws.Range("P108").Value = 3.69
$ws.Range("Q108").Value = "28/10/2023 14:16"
$ws.Range("R108").Value = 2.16
$ws.Range("S108").Value = "27/10/2023 02:42"
$ws.Range("T108").Value = 2.68
$ws.Range("U108").Value = "28/10/2023 14:15"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/slovakia/2-liga/trebisov-zilina/rkRSFJrA/"

# Row 109
$ws.Range("F109").Value = "Petrzalka"
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = "Komarno"
$ws.Range("I109").Value = 1
$ws.Range("J109").Value = 2.89
$ws.Range("K109").Value = "27/10/2023 23:42"
$ws.Range("L109").Value = 2.91
$ws.Range("M109").Value = "29/10/2023 10:02"
$ws.Range("N109").Value = 3.2
$ws.Range("O109").Value = "27/10/2023 23:42"
$ws.Range("P109").Value = 3.44
$ws.Range("Q109").Value = "29/10/2023 10:02"
$ws.Range("R109").Value = 2.2
$ws.Range("S109").Value = "27/10/2023 23:42"
$ws.Range("T109").Value = 2.29
$ws.Range("U109").Value = "29/10/2023 10:02"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/slovakia/2-liga/petrzalka-komarno/GdEF6F5d/"

# Row 110
$ws.Range("F110").Value = "Malzenice"
$ws.Range("G110").Value = 3
$ws.Range("H110").Value = "Spisska Nova Ves"
$ws.Range("I110").Value = 1
$ws.Range("J110").Value = 1.94
$ws.Range("K110").Value = "28/10/2023 08:13"
$ws.Range("L110").Value = 2.03
$ws.Range("M110").Value = "29/10/2023 10:21"
$ws.Range("N110").Value = 3.34
$ws.Range("O110").Value = "28/10/2023 08:13"
$ws.Range("P110").Value = 3.33
$ws.Range("Q110").Value = "29/10/2023 10:29"
$ws.Range("R110").Value = 3.32
$ws.Range("S110").Value = "28/10/2023 08:13"
$ws.Range("T110").Value = 3.59
$ws.Range("U110").Value = "29/10/2023 10:21"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-spisska-nova-ves/KUK68yyp/"

# --- Append two new rows (117, 118) at the bottom ---
# Copy the last existing data row's formatting (styles/number formats) down
# into the two new rows first, then overwrite every cell with its real value.
$ws.Range("A116:V116").Copy($ws.Range("A117:V117"))
$ws.Range("A116:V116").Copy($ws.Range("A118:V118"))

# Row 117
$ws.Range("A117").Value = 116
$ws.Range("B117").Value = "slovakia"
$ws.Range("C117").Value = "2-liga"
$ws.Range("D117").Value = "2023-2024"
$ws.Range("E117").Value = 45235.4375
$ws.Range("F117").Value = "Malzenice"
$ws.Range("G117").Value = 4
$ws.Range("H117").Value = "Samorin"
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2.06
$ws.Range("K117").Value = "05/11/2023 08:04"
$ws.Range("L117").Value = 2.06
$ws.Range("M117").Value = "05/11/2023 08:04"
$ws.Range("N117").Value = 3.52
$ws.Range("O117").Value = "05/11/2023 08:30"
$ws.Range("P117").Value = 3.52
$ws.Range("Q117").Value = "05/11/2023 08:30"
$ws.Range("R117").Value = 3.27
$ws.Range("S117").Value = "05/11/2023 08:04"
$ws.Range("T117").Value = 3.27
$ws.Range("U117").Value = "05/11/2023 08:04"
$ws.Range("V117").Value = "https://www.betexplorer.com/football/slovakia/2-liga/malzenice-samorin/4tZhKYZe/"

# Row 118
$ws.Range("A118").Value = 117
$ws.Range("B118").Value = "slovakia"
$ws.Range("C118").Value = "2-liga"
$ws.Range("D118").Value = "2023-2024"
$ws.Range("E118").Value = 45235.4375
$ws.Range("F118").Value = "Slovan Bratislava B"
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = "Puchov"
$ws.Range("I118").Value = 4
$ws.Range("J118").Value = 2.59
$ws.Range("K118").Value = "03/11/2023 22:44"
$ws.Range("L118").Value = 2.59
$ws.Range("M118").Value = "03/11/2023 22:44"
$ws.Range("N118").Value = 3.43
$ws.Range("O118").Value = "05/11/2023 08:33"
$ws.Range("P118").Value = 3.43
$ws.Range("Q118").Value = "05/11/2023 08:33"
$ws.Range("R118").Value = 2.4
$ws.Range("S118").Value = "03/11/2023 22:44"
$ws.Range("T118").Value = 2.4
$ws.Range("U118").Value = "03/11/2023 22:44"
$ws.Range("V118").Value = "https://www.betexplorer.com/football/slovakia/2-liga/slovan-bratislava-msk-puchov/AJZlLEKk/"
